# Refresh Leve profit-tracking figures (columns H-N) across multiple sheets
# with updated market-board pricing, as pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 52: Your Courtesy Wake-up Call / Smelling Salts
$ws.Range("H52").Value = 598.125
$ws.Range("I52").Value = 2000
$ws.Range("J52").Value = 397.85715
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 1193.57145
$ws.Range("M52").Value = -5840
$ws.Range("N52").Value = -1513.57145
# Row 58: A Matter of Vital Importance / Mega-Potion of Vitality
$ws.Range("H58").Value = 1125.2
$ws.Range("I58").Value = 90.40000000000001
$ws.Range("J58").Value = 2160
$ws.Range("K58").Value = 271.2
$ws.Range("L58").Value = 6480
$ws.Range("M58").Value = -121.2
$ws.Range("N58").Value = -6780
# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 6879.029
$ws.Range("I69").Value = 13000
$ws.Range("K69").Value = 39000
$ws.Range("M69").Value = -38126
# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 6879.029
$ws.Range("I72").Value = 13000
$ws.Range("K72").Value = 117000
$ws.Range("M72").Value = -112632
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 3701
$ws.Range("I76").Value = 3701
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3701
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3386
$ws.Range("N76").ClearContents()
# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 3701
$ws.Range("I79").Value = 3701
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3701
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2609
$ws.Range("N79").ClearContents()
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1612.238
$ws.Range("J138").Value = 1422.5
$ws.Range("L138").Value = 4267.5
$ws.Range("N138").Value = -14547.5

$ws = $wb.Worksheets.Item("ARM")
# Row 22: Kiss the Pan (Good-bye) / Initiate's Skillet
$ws.Range("H22").Value = 31966.5
$ws.Range("I22").Value = 31966.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 31966.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -31667.5
$ws.Range("N22").ClearContents()
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3711469.2
$ws.Range("I32").Value = 4389.1904
$ws.Range("K32").Value = 4389.1904
$ws.Range("M32").Value = -4102.1904
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3290.842
$ws.Range("I61").Value = 3189.8823
$ws.Range("J61").Value = 4149
$ws.Range("K61").Value = 3189.8823
$ws.Range("L61").Value = 4149
$ws.Range("M61").Value = -2977.8823
$ws.Range("N61").Value = -4573
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 4812.7
$ws.Range("I74").Value = 4105.8823
$ws.Range("K74").Value = 4105.8823
$ws.Range("M74").Value = -3231.8823
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 4812.7
$ws.Range("I77").Value = 4105.8823
$ws.Range("K77").Value = 20529.4115
$ws.Range("M77").Value = -16161.4115
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2899.8333
$ws.Range("I122").Value = 2899.8333
$ws.Range("K122").Value = 8699.499899999999
$ws.Range("M122").Value = -6249.499899999999
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 4159.2856
$ws.Range("I132").Value = 2507.3333
$ws.Range("J132").Value = 5398.25
$ws.Range("K132").Value = 7521.999899999999
$ws.Range("L132").Value = 16194.75
$ws.Range("M132").Value = -4991.999899999999
$ws.Range("N132").Value = -21254.75
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3290.842
$ws.Range("I136").Value = 3189.8823
$ws.Range("J136").Value = 4149
$ws.Range("K136").Value = 9569.6469
$ws.Range("L136").Value = 12447
$ws.Range("M136").Value = -7019.6469
$ws.Range("N136").Value = -17547

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 558
$ws.Range("I22").Value = 482.66666
$ws.Range("J22").Value = 633.3333
$ws.Range("K22").Value = 482.66666
$ws.Range("L22").Value = 633.3333
$ws.Range("M22").Value = -309.66666
$ws.Range("N22").Value = -979.3333
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2229.875
$ws.Range("I134").Value = 2229.875
$ws.Range("K134").Value = 6689.625
$ws.Range("M134").Value = -4154.625

$ws = $wb.Worksheets.Item("CRP")
# Row 41: The Lone Bowman / Oak Longbow
$ws.Range("H41").Value = 599430.25
$ws.Range("J41").Value = 883103.7
$ws.Range("L41").Value = 883103.7
$ws.Range("N41").Value = -883959.7
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 7512.6
$ws.Range("I58").Value = 1599
$ws.Range("K58").Value = 1599
$ws.Range("M58").Value = -1396
# Row 59: Bow Down to Magic / Crab Bow
$ws.Range("H59").Value = 36197.25
$ws.Range("J59").Value = 41082.57
$ws.Range("L59").Value = 41082.57
$ws.Range("N59").Value = -43372.57
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 2351.6667
$ws.Range("I62").Value = 2027.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2027.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1403.5
$ws.Range("N62").Value = -4248
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2351.6667
$ws.Range("I65").Value = 2027.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10137.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -7017.5
$ws.Range("N65").Value = -21240
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 4414.5
$ws.Range("I132").Value = 3550.5833
$ws.Range("K132").Value = 10651.7499
$ws.Range("M132").Value = -8121.749899999999
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1470.0526
$ws.Range("I134").Value = 1470.0526
$ws.Range("K134").Value = 4410.1578
$ws.Range("M134").Value = -1875.1578
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 7512.6
$ws.Range("I136").Value = 1599
$ws.Range("K136").Value = 4797
$ws.Range("M136").Value = -2247

$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand / Raisins
$ws.Range("H7").Value = 23.666666
$ws.Range("I7").Value = 23.666666
$ws.Range("K7").Value = 70.99999800000001
$ws.Range("M7").Value = 41.00000199999999
# Row 33: Cooking with Gas / Chicken Stock
$ws.Range("H33").Value = 101.666664
$ws.Range("I33").Value = 89.333336
$ws.Range("K33").Value = 536.000016
$ws.Range("M33").Value = -253.000016
# Row 92: Oh No Udon / Gyr Abanian Flour
$ws.Range("H92").Value = 2171.2222
$ws.Range("I92").Value = 780
$ws.Range("J92").Value = 3284.2
$ws.Range("K92").Value = 2340
$ws.Range("L92").Value = 9852.599999999999
$ws.Range("M92").Value = -1092
$ws.Range("N92").Value = -12348.6
# Row 109: Cure for What Ails / Purple Carrot Juice
$ws.Range("H109").Value = 2402.2222
$ws.Range("I109").Value = 1787.2
$ws.Range("J109").Value = 3171
$ws.Range("K109").Value = 5361.6
$ws.Range("L109").Value = 9513
$ws.Range("M109").Value = -4321.6
$ws.Range("N109").Value = -11593
# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = 1940
# Row 130: Blast from the Pasta / The Noodles of Elpis
$ws.Range("H130").Value = 3250
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 3250
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 9750
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -19790
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 3461.1667
$ws.Range("J131").Value = 5259.8
$ws.Range("L131").Value = 15779.4
$ws.Range("N131").Value = -25859.4
# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 11118.6
$ws.Range("I134").Value = 4865
$ws.Range("K134").Value = 14595
$ws.Range("M134").Value = -9525
# Row 136: Simple Is Hardest / Spaghetti al Olio e Peperoncino
$ws.Range("H136").Value = 9999.75
$ws.Range("I136").Value = 7500
$ws.Range("J136").Value = 10833
$ws.Range("K136").Value = 22500
$ws.Range("L136").Value = 32499
$ws.Range("M136").Value = -17400
$ws.Range("N136").Value = -42699
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 5533
$ws.Range("J138").Value = 6799.5
$ws.Range("L138").Value = 20398.5
$ws.Range("N138").Value = -30678.5
# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 3541.5
$ws.Range("I139").Value = 2243.2
$ws.Range("J139").Value = 10033
$ws.Range("K139").Value = 6729.599999999999
$ws.Range("L139").Value = 30099
$ws.Range("M139").Value = -1589.599999999999
$ws.Range("N139").Value = -40379

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 42765.223
$ws.Range("I132").Value = 57617.156
$ws.Range("J132").Value = 7491.875
$ws.Range("K132").Value = 172851.468
$ws.Range("L132").Value = 22475.625
$ws.Range("M132").Value = -170321.468
$ws.Range("N132").Value = -27535.625

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 802.625
$ws.Range("I16").Value = 745.8570999999999
$ws.Range("K16").Value = 745.8570999999999
$ws.Range("M16").Value = -575.8570999999999
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 5708.5
$ws.Range("I122").Value = 5708.5
$ws.Range("K122").Value = 17125.5
$ws.Range("M122").Value = -14675.5
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2730.4
$ws.Range("I132").Value = 2688.8235
$ws.Range("J132").Value = 2966
$ws.Range("K132").Value = 8066.470499999999
$ws.Range("L132").Value = 8898
$ws.Range("M132").Value = -5536.470499999999
$ws.Range("N132").Value = -13958

$ws = $wb.Worksheets.Item("WVR")
# Row 9: A Taste for Dalmaticae / Amateur's Dalmatica
$ws.Range("H9").Value = 700
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
# Row 62: Pride Up in Smoke / Rainbow Cloth
$ws.Range("H62").Value = 11332.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11332.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11332.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -12580.5
# Row 65: Desperate for Diversionaries (L) / Rainbow Cloth
$ws.Range("H65").Value = 11332.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11332.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 56662.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -62902.5
# Row 95: Duress Rehearsal / Ruby Cotton Fingerless Gloves of Casting
$ws.Range("H95").Value = 43579.6
$ws.Range("J95").Value = 43579.6
$ws.Range("L95").Value = 43579.6
$ws.Range("N95").Value = -49071.6
# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 5002
$ws.Range("I100").Value = 5002
$ws.Range("K100").Value = 10004
$ws.Range("M100").Value = -9463
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 962
$ws.Range("I122").Value = 962
$ws.Range("K122").Value = 2886
$ws.Range("M122").Value = -436
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2679
$ws.Range("I132").Value = 2376.8572
$ws.Range("K132").Value = 7130.571599999999
$ws.Range("M132").Value = -4600.571599999999
